$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$ws.Range("A7").Value = "Mike"
$ws.Range("B7").Value = "Jack"
$ws.Range("C7").Value = "Mii"
$ws.Range("D7").Value = "02:50.610"
$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E7").Value = 45982

$ws.Range("C8").Select()
